$d = $word.ActiveDocument

# 1) //span[contains(text(),'Select a Country')]/parent::a
#    -> //span[contains(text(),'Select a Country')]//ancestor::div/preceding-sibling::select
$rng = $d.Content
$found = $rng.Find.Execute("//span[contains(text(),'Select a Country')]/parent::a")
if ($found) {
    $rng.Text = "//span[contains(text(),'Select a Country')]//ancestor::div/preceding-sibling::select"
}

# 2) //span[contains(text(),'Select a Category')]/parent::a
#    -> //span[contains(text(),'Select a Category')]//ancestor::div/preceding-sibling::select
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("//span[contains(text(),'Select a Category')]/parent::a")
if ($found2) {
    $rng2.Text = "//span[contains(text(),'Select a Category')]//ancestor::div/preceding-sibling::select"
}

# 3) Split the run holding //div[contains(text(),'Select State')]/parent::div//input
#    right after "//div[contains" and drop the "_GoBack" bookmark there
#    (it moves from the trailing empty paragraph to this split point).
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("//div[contains(text(),'Select State')]/parent::div//input")
if ($found3) {
    $splitPos = $rng3.Start + 14
    $bmRange = $d.Range($splitPos, $splitPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# 4) styles.xml: mark FollowedHyperlink / Normal Table as Quick Styles (w:qFormat)
$d.Styles("FollowedHyperlink").QuickStyle = $true
$d.Styles("Normal Table").QuickStyle = $true
